$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.692.62'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '3.370.78'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.88'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.64'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.369.37'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.58'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('E11').Value = '  -3.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.379'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('D13').Value = '3.943.77'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.78'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').Value = '3.373.06'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000169'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -3.73%  '
$ws.Range('D18').Value = '60.859.31'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.79'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.68'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -3.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.23'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -2.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '372.29'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').Value = '3.509.86'
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('E24').Value = '  -2.40%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.81'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.176'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +10.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.62'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -4.04%  '
$ws.Range('E30').Value = '  -0.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.28'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -3.53%  '
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('E33').Value = '  -2.02%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.22'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.10'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -4.46%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.76'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -0.93%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.53'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -2.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.55'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0754'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -3.41%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.771'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -1.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '24.92'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('E44').Value = '  -2.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.29'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('E46').Value = '  -5.75%  '
$ws.Range('D47').Value = '2.523.66'
$ws.Range('E47').Value = '  +7.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.19'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +2.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.74'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -1.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.41'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +4.25%  '
$ws.Range('E51').Value = '  -1.52%  '
